$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the I (TIPO DE DATO) / J (FORMATO) columns for rows 34-51: the
# stale "NA" placeholder is replaced with the real image format info.
for ($r = 34; $r -le 51; $r++) {
    $ws.Range("I$r").Value = "Imagen"
    $ws.Range("J$r").Value = "png"
}

# Row 51 was missing the bottom border that the rest of the table block
# (rows 34-50) already has on columns I/J; restore it so the cell picks
# up the same boxed style as its neighbours.
$r51 = $ws.Range("I51:J51")
$r51.Borders.Item(9).LineStyle = 1
$r51.Borders.Item(9).Weight = 2
$r51.Borders.Item(9).ColorIndex = 1

# The leftover _FilterDatabase defined name still pointed at the old
# single-row range; repoint it at the full data block (A1:P118).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$P`$118"
    }
}

# Move the frozen-pane selection to where review work left off.
$ws.Range("D54").Select()
